$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: "Spring" part with a McMaster-Carr hyperlink, matching the
# pattern used by the existing rows 2-4 (plain text label in column A,
# hyperlinked URL with the "Hyperlink" style in column B).
$ws.Range("A5").Value = "Spring"
$ws.Range("B5").Value = "https://www.mcmaster.com/9657K248/"

$ws.Hyperlinks.Add($ws.Range("B5"), "https://www.mcmaster.com/9657K248/")
$ws.Range("B5").Style = "Hyperlink"

# Move the active selection to D10, as in the target workbook.
[void]$ws.Range("D10").Select()
